# TC04_Canine_Filter_FileFormat-pdf.xlsx edit
#
# The "Cases" Neo4j query stored in cell B2 of the "startup" sheet contained
# a trailing `coalesce(co.cohort_description, '') AS `Cohort`` column that is
# being removed (the Cohort information is no longer returned by that
# query). B3/B4 keep their existing values unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Cases query text in B2: drop the trailing Cohort column ---
$b2 = $ws.Range("B2").Value()
$lines = $b2 -split "`n"

# Remove the last line (the Cohort coalesce) ...
$newLines = $lines[0..($lines.Length - 2)]
# ... and strip the now-trailing comma from the line that used to precede it.
$lastIdx = $newLines.Length - 1
$newLines[$lastIdx] = $newLines[$lastIdx].TrimEnd(",")

$newB2 = [string]::Join("`n", $newLines)
$ws.Range("B2").Value = $newB2

# The author's selection ended up on B2 after making this edit.
[void]$ws.Range("B2").Select()

# Row heights shrink slightly because the wrapped text now renders with one
# fewer line / a tighter line metric.
$ws.Rows(2).RowHeight = 259.2
$ws.Rows(3).RowHeight = 288
$ws.Rows(4).RowHeight = 259.2
